$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.ClearFormats()
}

Set-TextValue "D2" "56.877.32"
Set-TextValue "D3" "2.968.24"
Set-TextValue "E3" "  -1.57%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "499.61"
Set-TextValue "E5" "  -3.37%  "
Set-TextValue "D6" "137.86"
Set-TextValue "E6" "  -1.48%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.04%  "
Set-TextValue "D8" "0.428"
Set-TextValue "E8" "  -2.15%  "
Set-TextValue "D9" "7.32"
Set-TextValue "E9" "  -3.43%  "
Set-TextValue "E10" "  -2.62%  "
Set-TextValue "D11" "0.357"
Set-TextValue "E11" "  -0.60%  "
Set-TextValue "D12" "3.471.24"
Set-TextValue "E12" "  -1.60%  "
Set-TextValue "E13" "  -1.80%  "
Set-TextValue "D14" "25.84"
Set-TextValue "E14" "  -0.22%  "
Set-TextValue "D15" "0.0000158"
Set-TextValue "E15" "  -0.50%  "
Set-TextValue "D16" "56.943.47"
Set-TextValue "E16" "  -0.20%  "
Set-TextValue "D17" "6.07"
Set-TextValue "E17" "  +1.15%  "
Set-TextValue "D18" "2.967.18"
Set-TextValue "E18" "  -1.54%  "
Set-TextValue "D19" "12.60"
Set-TextValue "E19" "  -0.53%  "
Set-TextValue "D20" "7.81"
Set-TextValue "E20" "  -1.36%  "
Set-TextValue "D21" "318.93"
Set-TextValue "E21" "  -3.28%  "
Set-TextValue "E22" "  -0.16%  "
Set-TextValue "D23" "5.65"
Set-TextValue "E23" "  -0.76%  "
Set-TextValue "D24" "0.484"
Set-TextValue "E24" "  -0.34%  "
Set-TextValue "D25" "63.12"
Set-TextValue "E25" "  -1.14%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.12%  "
Set-TextValue "D27" "0.163"
Set-TextValue "E27" "  -5.32%  "
Set-TextValue "D28" "0.0₃0892"
Set-TextValue "E28" "  -3.15%  "
Set-TextValue "E29" "  -3.03%  "
Set-TextValue "D30" "7.08"
Set-TextValue "E30" "  -1.93%  "
Set-TextValue "D31" "1.76"
Set-TextValue "E31" "  -3.29%  "
Set-TextValue "D32" "1.16"
Set-TextValue "E32" "  -5.55%  "
Set-TextValue "D33" "20.09"
Set-TextValue "E33" "  -3.15%  "
Set-TextValue "D34" "155.12"
Set-TextValue "E34" "  -1.60%  "
Set-TextValue "D35" "4.61"
Set-TextValue "E35" "  -0.46%  "
Set-TextValue "D36" "5.73"
Set-TextValue "E36" "  -0.43%  "
Set-TextValue "D37" "1.24"
Set-TextValue "E37" "  -3.25%  "
Set-TextValue "D38" "24.06"
Set-TextValue "E38" "  -0.69%  "
Set-TextValue "D39" "0.0666"
Set-TextValue "E39" "  -2.37%  "
Set-TextValue "B40" "OKB"
Set-TextValue "C40" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D40" "37.59"
Set-TextValue "E40" "  +0.79%  "
Set-TextValue "B41" "RenzoRestakedETH"
Set-TextValue "C41" "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue "D41" "2.996.00"
Set-TextValue "E41" "  -1.64%  "
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  +0.13%  "
Set-TextValue "D43" "3.72"
Set-TextValue "E43" "  -0.26%  "
Set-TextValue "D44" "0.638"
Set-TextValue "E44" "  -2.00%  "
Set-TextValue "D45" "2.201.82"
Set-TextValue "E45" "  -4.15%  "
Set-TextValue "D46" "1.38"
Set-TextValue "E46" "  -3.65%  "
Set-TextValue "D47" "0.943"
Set-TextValue "E47" "  -6.34%  "
Set-TextValue "D48" "5.92"
Set-TextValue "E48" "  +0.64%  "
Set-TextValue "E49" "  -2.99%  "
Set-TextValue "D50" "19.23"
Set-TextValue "E51" "  -10.31%  "
